$wb = $excel.ActiveWorkbook

# --- Sheet1 (pcroprep) ---
$ws1 = $wb.Worksheets.Item("pcroprep")

# Row 35 (rice_wheat / veges): remove D35 and F35, update G35
$ws1.Range("D35").ClearContents()
$ws1.Range("F35").ClearContents()
$ws1.Range("G35").Value = -241.4

# Row 39 (rice_wheat / total): update D39, F39, G39
$ws1.Range("D39").Value = 940.06243687239601
$ws1.Range("F39").Value = 284.77792550326666
$ws1.Range("G39").Value = -521.82207449673342

# --- Sheet4 (pdietrep) ---
$ws4 = $wb.Worksheets.Item("pdietrep")

$ws4.Range("E6").Value = 1214.8943640250632
$ws4.Range("F6").Value = -942.89115018061966
$ws4.Range("G6").Value = 56.302832511705233

$ws4.Range("E7").Value = 41.227792490760059
$ws4.Range("F7").Value = -28.900398385590499
$ws4.Range("G7").Value = 58.7891858831102

$ws4.Range("E8").Value = 14.785465764643714
$ws4.Range("F8").Value = -49.948099661526754
$ws4.Range("G8").Value = 22.840493440001758

$ws4.Range("E9").Value = 215.13404923941323
$ws4.Range("F9").Value = -108.53377789143912
$ws4.Range("G9").Value = 66.467542092912083

# --- Sheet5 (pradar) ---
$ws5 = $wb.Worksheets.Item("pradar")

# Row 15 (rice_wheat / veget): remove D15 and E15, update F15
$ws5.Range("D15").ClearContents()
$ws5.Range("E15").ClearContents()
$ws5.Range("F15").Value = -241.4

# --- Sheet6 (plandrep) ---
$ws6 = $wb.Worksheets.Item("plandrep")

# Row 11: remove S11
$ws6.Range("S11").ClearContents()

# --- Sheet7 (plaborrep) ---
$ws7 = $wb.Worksheets.Item("plaborrep")

# Row 3: remove R3, update AF3
$ws7.Range("R3").ClearContents()
$ws7.Range("AF3").Value = 0.91221062026425437

# --- Sheet8 (pfertrep) ---
$ws8 = $wb.Worksheets.Item("pfertrep")

# Row 5: remove S5, update Z5
$ws8.Range("S5").ClearContents()
$ws8.Range("Z5").Value = 291586.2928207317

# Row 6: remove S6, update Z6
$ws8.Range("S6").ClearContents()
$ws8.Range("Z6").Value = 341415.28265000001

# Row 7: remove S7, update Z7
$ws8.Range("S7").ClearContents()
$ws8.Range("Z7").Value = 323916.94579756091

$wb.Save()
